# Auto-generated Excel COM-interop script to apply cryptos.xlsx data refresh
# (GitHub Actions style update of prices/volumes, and an Aptos/WrappedeETH row swap)
#
# All target cells in this sheet are stored as text (t="inlineStr") even when
# their contents look numeric (prices/percentages). Excel auto-converts plain
# numeric-looking strings to real numbers on assignment, which would silently
# change values like '0.200' -> 0.2 or '237.29' -> 237.28999999999999. To avoid
# this, we force the NumberFormat of those cells to Text ("@") before writing
# the new value so the literal string is preserved exactly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '95.194.67'
$ws.Cells.Item(2, 5).Value = '  +1.89%  '
$ws.Cells.Item(3, 4).Value = '3.614.74'
$ws.Cells.Item(3, 5).Value = '  +5.96%  '
$ws.Cells.Item(4, 5).Value = '  +0.04%  '
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '237.29'
$ws.Cells.Item(5, 5).Value = '  +1.56%  '
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = '658.64'
$ws.Cells.Item(6, 5).Value = '  +6.17%  '
$ws.Cells.Item(7, 5).Value = '  +0.31%  '
$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = '0.404'
$ws.Cells.Item(8, 5).Value = '  +3.28%  '
$ws.Cells.Item(9, 5).Value = '  -0.02%  '
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = '0.991'
$ws.Cells.Item(10, 5).Value = '  -0.19%  '
$ws.Cells.Item(11, 4).Value = '3.614.27'
$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = '42.42'
$ws.Cells.Item(12, 5).Value = '  -2.52%  '
$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = '0.200'
$ws.Cells.Item(13, 5).Value = '  +0.61%  '
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = '6.29'
$ws.Cells.Item(14, 5).Value = '  -0.01%  '
$ws.Cells.Item(15, 4).Value = '4.290.47'
$ws.Cells.Item(15, 5).Value = '  +5.63%  '
$ws.Cells.Item(16, 4).Value = '95.404.10'
$ws.Cells.Item(16, 5).Value = '  +2.26%  '
$ws.Cells.Item(17, 5).Value = '  +2.80%  '
$ws.Cells.Item(18, 4).Value = '3.615.95'
$ws.Cells.Item(18, 5).Value = '  +6.11%  '
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = '12.99'
$ws.Cells.Item(19, 5).Value = '  +11.59%  '
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = '7.75'
$ws.Cells.Item(20, 5).Value = '  -6.72%  '
$ws.Cells.Item(21, 5).Value = '  -0.34%  '
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = '3.57'
$ws.Cells.Item(22, 5).Value = '  +4.95%  '
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = '0.480'
$ws.Cells.Item(23, 5).Value = '  -6.18%  '
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = '505.19'
$ws.Cells.Item(24, 5).Value = '  +1.23%  '
$ws.Cells.Item(25, 5).Value = '  +7.27%  '
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = '6.60'
$ws.Cells.Item(26, 5).Value = '  -3.25%  '
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = '95.74'
$ws.Cells.Item(27, 5).Value = '  +8.67%  '
$ws.Cells.Item(28, 2).Value = 'WrappedeETH'
$ws.Cells.Item(28, 3).Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Cells.Item(28, 4).Value = '3.815.77'
$ws.Cells.Item(28, 5).Value = '  +6.24%  '
$ws.Cells.Item(29, 2).Value = 'Aptos'
$ws.Cells.Item(29, 3).Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = '12.55'
$ws.Cells.Item(29, 5).Value = '  +4.20%  '
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = '3.16'
$ws.Cells.Item(30, 5).Value = '  +16.19%  '
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = '11.32'
$ws.Cells.Item(31, 5).Value = '  -0.36%  '
$ws.Cells.Item(32, 5).Value = '  -0.23%  '
$ws.Cells.Item(33, 5).Value = '  -1.97%  '
$ws.Cells.Item(34, 5).Value = '  -0.13%  '
$ws.Cells.Item(35, 5).Value = '  +0.62%  '
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = '31.93'
$ws.Cells.Item(36, 5).Value = '  +10.44%  '
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = '0.557'
$ws.Cells.Item(37, 5).Value = '  +0.83%  '
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = '573.56'
$ws.Cells.Item(38, 5).Value = '  +2.85%  '
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = '8.11'
$ws.Cells.Item(39, 5).Value = '  +8.11%  '
$ws.Cells.Item(40, 5).Value = '  +5.06%  '
$ws.Cells.Item(41, 5).Value = '  +0.03%  '
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = '0.921'
$ws.Cells.Item(42, 5).Value = '  +3.07%  '
$ws.Cells.Item(43, 5).Value = '  -0.30%  '
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = '34.94'
$ws.Cells.Item(44, 5).Value = '  +46.56%  '
$ws.Cells.Item(45, 5).Value = '  +1.47%  '
$ws.Cells.Item(46, 5).Value = '  -0.23%  '
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = '5.64'
$ws.Cells.Item(47, 5).Value = '  +3.08%  '
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = '2.22'
$ws.Cells.Item(48, 5).Value = '  +5.45%  '
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = '0.0411'
$ws.Cells.Item(49, 5).Value = '  -1.75%  '
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = '3.49'
$ws.Cells.Item(50, 5).Value = '  -3.94%  '
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = '53.59'
$ws.Cells.Item(51, 5).Value = '  +0.80%  '
